# Manual Operario.docx - "cambio de fechas a todas las actividades de la
# ultima entrega": the cover-page delivery date moves from 25/10/2019 to
# 1/11/2019, which is what Word's edit-tracking leaves behind as a shuffled
# "_GoBack" bookmark (and a consequent renumbering of the other bookmarks).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: merge the two runs that make up the "...icono de una "i"." text.
# In the original document this sentence was split into two runs with a
# stray "_GoBack" bookmark right after it; retyping the whole sentence in
# one shot (like a user re-typing it) naturally collapses it back into a
# single run with the bookmark gone from here - it will reappear next to
# the new date below.
# ---------------------------------------------------------------------
$startRng = $d.Content.Duplicate
$startRng.Find.Execute("Para ingresar a la información acerca de nuestra empresa")
$mergeStart = $startRng.Start

$endRng = $d.Content.Duplicate
$endRng.Find.Execute("una “i”.")
$mergeEnd = $endRng.End

$mergeRange = $d.Range($mergeStart, $mergeEnd)
$mergeRange.Select()
$word.Selection.TypeText("Para ingresar a la información acerca de nuestra empresa y aplicación se debe hacer click en el botón ubicado en la esquina superior derecha de la aplicación, el cual posee un icono de una “i”.")

# ---------------------------------------------------------------------
# Step 2: change the cover date from 25/10/2019 to 1/11/2019. Word leaves
# the cursor (and hence the "_GoBack" bookmark) right after the "1/11" that
# was just typed, splitting that run away from the trailing "/2019".
# ---------------------------------------------------------------------
$dateRng = $d.Content.Duplicate
$dateRng.Find.Execute("25/10/2019")
$dateStart = $dateRng.Start
$dateEnd = $dateRng.End

$dateTarget = $d.Range($dateStart, $dateEnd)
$dateTarget.Text = "1/11/2019"

# Force the run boundary right after "Tercera entrega " too, so the
# surviving text is split into three runs, matching a real edit session.
$splitPoint = $d.Range($dateStart, $dateStart)
$d.Bookmarks.Add("_TmpSplit", $splitPoint)

# Re-adding "_GoBack" here moves it from its old location (Step 1's
# sentence) to right after "1/11", which also splits "/2019" into its own
# run.
$goBackPos = $dateStart + 4
$goBackPoint = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackPoint)

$d.Bookmarks("_TmpSplit").Delete()
